$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a handful of existing one-handed entries (error-check fixes) ---
$ws.Range("E5").Value = "-20,0,1,1,1,1,1,1"
$ws.Range("D12").Value = "-9,4,8,8"
$ws.Range("E13").Value = "R,14,6:L,14,2,7,5,6,5"
$ws.Range("E14").Value = "R,14,6:L,14,3,6,6,8"
$ws.Range("E18").Value = "R,15,6:L,15,2,5,7,5,5"

# --- New hangboard entries: rows 19-28 ---
$ws.Cells.Item(19,1).Value = "1 Mar 2023"
$ws.Cells.Item(19,2).Value = "0,6"
$ws.Cells.Item(19,3).Value = "-4.5,6"
$ws.Cells.Item(19,4).Value = "-8.5,4,7,6"
$ws.Cells.Item(19,5).Value = "L,14,6:R,14,6"
$ws.Cells.Item(19,6).Value = "-15,6"
$ws.Cells.Item(19,7).Value = "-5.5,6"
$ws.Cells.Item(19,8).Value = "L,12,5,5:R,12,5,9"
$ws.Cells.Item(19,9).Value = "-17,6"
$ws.Cells.Item(19,10).Value = "-14.5,6"
$ws.Cells.Item(20,1).Value = "3 Mar 2023"
$ws.Cells.Item(20,2).Value = "0,6"
$ws.Cells.Item(20,3).Value = "-4,6"
$ws.Cells.Item(20,4).Value = "-8.5,5,7"
$ws.Cells.Item(20,5).Value = "L,14.5,3,7,8,4:R,14.5,6"
$ws.Cells.Item(20,6).Value = "-14.5,6"
$ws.Cells.Item(20,7).Value = "-5,6"
$ws.Cells.Item(20,8).Value = "L,12.5,4,7,5:R,12.5,5,7"
$ws.Cells.Item(20,9).Value = "-16.5,5,7"
$ws.Cells.Item(20,10).Value = "-14,4,7,7"
$ws.Cells.Item(21,1).Value = "8 Mar 2023"
$ws.Cells.Item(21,2).Value = "0,6"
$ws.Cells.Item(21,3).Value = "-3.5,4,9,5"
$ws.Cells.Item(21,4).Value = "-8.5,3,5,8,5"
$ws.Cells.Item(21,5).Value = "L,14.5,6:R,14.5,6"
$ws.Cells.Item(21,6).Value = "-14,6"
$ws.Cells.Item(21,7).Value = "-4.5,6"
$ws.Cells.Item(21,8).Value = "L,12.5,4,7,4:R,12.5,4,4,2"
$ws.Cells.Item(21,9).Value = "-16.5,5,6"
$ws.Cells.Item(21,10).Value = "-14,5,6"
$ws.Cells.Item(22,1).Value = "13 Mar 2023"
$ws.Cells.Item(22,2).Value = "0,6"
$ws.Cells.Item(22,3).Value = "-3.5,5,5"
$ws.Cells.Item(22,4).Value = "-8.5,5,4"
$ws.Cells.Item(22,5).Value = "L,14.5,3,8,5,6:R,14.5,5,9"
$ws.Cells.Item(22,6).Value = "-13.5,6"
$ws.Cells.Item(22,7).Value = "-4,5,9"
$ws.Cells.Item(22,8).Value = "L,12.5,2,5,5,5,4:R,12.5,4,8,5"
$ws.Cells.Item(22,9).Value = "-16.5,6"
$ws.Cells.Item(22,10).Value = "-14,6"
$ws.Cells.Item(23,1).Value = "22 Mar 2023"
$ws.Cells.Item(23,2).Value = "0,6"
$ws.Cells.Item(23,3).Value = "-3.5,5,4"
$ws.Cells.Item(23,4).Value = "-8.5,3,8,6,5"
$ws.Cells.Item(23,5).Value = "R,14.5,5,9:L,14.5,4,6,4"
$ws.Cells.Item(23,6).Value = "-13,5,8"
$ws.Cells.Item(23,7).Value = "-4,3,8,8,7"
$ws.Cells.Item(23,8).Value = "R,12.5,4,6,6:L,12.5,3,7,6,5"
$ws.Cells.Item(23,9).Value = "-16,3,9,8,6"
$ws.Cells.Item(23,10).Value = "-13.5,4,6,8"
$ws.Cells.Item(24,1).Value = "5 Apr 2023"
$ws.Cells.Item(24,2).Value = "0,6"
$ws.Cells.Item(24,3).Value = "-3.5,4,8,5"
$ws.Cells.Item(24,4).Value = "-8.5,5,7"
$ws.Cells.Item(24,5).Value = "R,14.5,6:L,14.5,4,8,5"
$ws.Cells.Item(24,6).Value = "-12.5,6"
$ws.Cells.Item(24,7).Value = "-4,5,7"
$ws.Cells.Item(24,8).Value = "R,12.5,5,9:L,12.5,2,7,7,7,6"
$ws.Cells.Item(24,9).Value = "-16,4,8,6"
$ws.Cells.Item(24,10).Value = "-13.5,5,8"
$ws.Cells.Item(25,1).Value = "14 Apr 2023"
$ws.Cells.Item(25,2).Value = "0,6"
$ws.Cells.Item(25,3).Value = "-3.5,5,9"
$ws.Cells.Item(25,4).Value = "-8,3,7,6,6"
$ws.Cells.Item(25,5).Value = "R,15,6:L,15,4,7,7"
$ws.Cells.Item(25,6).Value = "-12,6"
$ws.Cells.Item(25,7).Value = "-3.5,5,9"
$ws.Cells.Item(25,8).Value = "R,13,4,6,6:L,13,4,9,6"
$ws.Cells.Item(25,9).Value = "-16,3,8,8,7"
$ws.Cells.Item(25,10).Value = "-13,4,4,8,7"
$ws.Cells.Item(26,1).Value = "20 May 2023"
$ws.Cells.Item(26,2).Value = "0,6"
$ws.Cells.Item(26,3).Value = "-3.5,6"
$ws.Cells.Item(26,4).Value = "-8,4,6,6"
$ws.Cells.Item(26,5).Value = "R,15,6:L,15,3,8,5,5"
$ws.Cells.Item(26,6).Value = "-11.5,6"
$ws.Cells.Item(26,7).Value = "-3.5,6"
$ws.Cells.Item(26,8).Value = "R,13,3,9,7,7:L,13,2,8,6,5,4"
$ws.Cells.Item(26,9).Value = "-16,3,7,7,7"
$ws.Cells.Item(26,10).Value = "-13,5,7"
$ws.Cells.Item(27,1).Value = "10 Jul 2023"
$ws.Cells.Item(27,2).Value = "0,6"
$ws.Cells.Item(27,3).Value = "-3.5,5,8"
$ws.Cells.Item(27,4).Value = "-8.5,5,5"
$ws.Cells.Item(27,5).Value = "L,14.5,2,4,7,6,6:R,14.5,6"
$ws.Cells.Item(27,6).Value = "-13,6"
$ws.Cells.Item(27,7).Value = "-4,6"
$ws.Cells.Item(27,8).Value = "L,12.5,2,6,6,3,3:R,12.5,4,5,6"
$ws.Cells.Item(27,9).Value = "-16,4,7,6"
$ws.Cells.Item(27,10).Value = "-13.5,5,6"
$ws.Cells.Item(28,1).Value = "19 Jul 2023"
$ws.Cells.Item(28,2).Value = "0,6"
$ws.Cells.Item(28,3).Value = "-3,4,9,5"
$ws.Cells.Item(28,4).Value = "-8,4,5,6"
$ws.Cells.Item(28,5).Value = "R,15,6:L,14.5,5,9"
$ws.Cells.Item(28,6).Value = "-12.5,5,5"
$ws.Cells.Item(28,7).Value = "-3.5,5,9"
$ws.Cells.Item(28,8).Value = "R,12.5,3,9,6,7:L,12.5,3,6,5,5"
$ws.Cells.Item(28,9).Value = "-16,2,8,7,7,6"
$ws.Cells.Item(28,10).Value = "-13,3,6,8,7"

# --- Column width adjustments (closest achievable via Excel pixel grid) ---
$ws.Range("E1:F1").ColumnWidth = 21.17
$ws.Range("H1").ColumnWidth = 23.83

# --- Update selection to match the last-edited cell ---
$ws.Range("H28").Select()
